# Genome_assembly_report.xlsx edit:
#   1. Bump the template Version metadata value from 0.0.1 to 1.0.0.
#   2. Add a new "Output [Data]" column to the annotationTable
#      (table grows from A1:P1 to A1:Q1).

$wb = $excel.ActiveWorkbook

# --- 1. Update Version metadata on the isa_template sheet (B4: 0.0.1 -> 1.0.0)
$wsMeta = $wb.Worksheets.Item("isa_template")
$wsMeta.Range("B4").Value = "1.0.0"

# --- 2. Add "Output [Data]" column to the annotationTable on "New Table" sheet
$wsTable = $wb.Worksheets.Item("New Table")
$tbl = $wsTable.ListObjects.Item(1)
$newCol = $tbl.ListColumns.Add()
$newCol.Range.Cells(1, 1).Value = "Output [Data]"
